$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: clear all the placeholder inline-string cells (A3:P3) ---
# ClearContents alone drops the row entirely from the XML; toggling
# Hidden keeps a bare <row r="3"/> element (matches the target diff).
$ws.Range("A3:P3").ClearContents()
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(3).Hidden = $false

# --- Row 4: drop the stray empty M4 cell ---
$ws.Cells.Item(4, 13).ClearContents()

# --- Row 5: new sales-order entry ---
# Columns that look numeric/date-like to Excel's type inference need to be
# forced to Text first so they round-trip as inline strings like their
# row-2/row-4 counterparts instead of turning into dates/numbers.
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "02/19/2024"
$ws.Cells.Item(5, 2).Value = "SO240219002"
$ws.Cells.Item(5, 3).Value = "Me"
$ws.Cells.Item(5, 4).Value = "5555555555"
$ws.Cells.Item(5, 5).Value = "Type O Negative"
$ws.Cells.Item(5, 6).Value = "Bloody Kisses"
$ws.Cells.Item(5, 7).Value = 5
$ws.Cells.Item(5, 8).Value = 49.99
$ws.Cells.Item(5, 9).Value = "AMS"
$ws.Cells.Item(5, 10).Value = "DVD"
$ws.Cells.Item(5, 11).Value = "MTP"
$ws.Cells.Item(5, 12).Value = "YES"
$ws.Cells.Item(5, 13).NumberFormat = "@"
$ws.Cells.Item(5, 13).Value = "59848"
$ws.Cells.Item(5, 14).Value = "gh"
$ws.Cells.Item(5, 15).Value = "hg"
$ws.Cells.Item(5, 16).NumberFormat = "@"
$ws.Cells.Item(5, 16).Value = "77777"
